$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 1029 (shifts existing rows 1029-1079 down to 1032-1082)
$ws.Range("A1029:A1031").EntireRow.Insert()

# Seed the 3 new rows with a copy of the row immediately below them (which used
# to be the first "Especial"/"Primera"/"Segunda" trio before the insert), so all
# the non-price columns (A-L, Q, R, T) start out correct.
$ws.Range("A1032:T1032").Copy()
$ws.Range("A1029:T1029").PasteSpecial()
$ws.Range("A1033:T1033").Copy()
$ws.Range("A1030:T1030").PasteSpecial()
$ws.Range("A1034:T1034").Copy()
$ws.Range("A1031:T1031").PasteSpecial()

# New week's data (date 2023-03-23 / serial 45008) for the inserted rows.
# Row 1029: "Especial"
$ws.Range("D1029").Value = 45008
$ws.Range("M1029").Value = 360
$ws.Range("N1029").Value = 13000
$ws.Range("O1029").Value = 14000
$ws.Range("P1029").Value = 13500
$ws.Range("S1029").Value = 1929

# Row 1030: "Primera"
$ws.Range("D1030").Value = 45008
$ws.Range("M1030").Value = 280
$ws.Range("N1030").Value = 11000
$ws.Range("O1030").Value = 12000
$ws.Range("P1030").Value = 11500
$ws.Range("S1030").Value = 1643

# Row 1031: "Segunda"
$ws.Range("D1031").Value = 45008
$ws.Range("M1031").Value = 240
$ws.Range("N1031").Value = 9000
$ws.Range("O1031").Value = 10000
$ws.Range("P1031").Value = 9500
$ws.Range("S1031").Value = 1357
